# Handles float input without breaking stuff
#
# This script reproduces, via the Excel COM object model, the corrections that
# were made to the generated marksheet worksheet ("quiz" / 1401EE65):
#   - the summary block (rows 10-12) is recomputed with the correct
#     right/wrong/not-attempted counts, marking scheme and total, and the
#     row-label cells (A10/A11/A12) are given the same bold "title" style
#     used by the other row labels;
#   - the "Marking -> Wrong" cell (C11) is stored as a real number (-1)
#     instead of text;
#   - the duplicated 2nd/3rd "Student Ans/Correct Ans" blocks (columns D:E
#     and G:H) that were left over from a bug are removed, since the sheet
#     only has a single, 28 question, answer key (columns A:B);
#   - the "Student Ans" column (A) is filled in for every question that
#     previously had no recorded answer, using the green "correctStyle" for
#     answers that match the key and the red "incorrectStyle" for answers
#     that don't.
#
# xlPasteFormats is used (via Copy + PasteSpecial) whenever a cell's style
# needs to change, since that reuses the workbook's existing named cell
# styles/xf records instead of generating new, duplicate ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1. Remove the redundant 2nd ("Student Ans"/"Correct Ans" in D:E, for rows
#    19 and below) and 3rd (G:H) answer-key blocks entirely.
# ---------------------------------------------------------------------------
$ws.Range("G15:H21").Clear()
$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------------------
# 2. Give the row labels in the summary block the same "mtitleStyle" used
#    by the other labels on the sheet (copy format from A9, which already
#    uses it), then fix up the summary numbers themselves.
# ---------------------------------------------------------------------------
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial($xlPasteFormats)
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial($xlPasteFormats)
$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial($xlPasteFormats)

$ws.Range("A10").Value = "No."
$ws.Range("B10").Value = 17
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = 28

$ws.Range("A11").Value = "Marking"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("A12").Value = "Total"
$ws.Range("B12").Value = 68
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "66/112"

# ---------------------------------------------------------------------------
# 3. Fill in the previously-blank "Student Ans" cells (column A, rows 16-40)
#    with the recovered answers, colouring them with the worksheet's
#    existing "correctStyle" (green, matches the key) or "incorrectStyle"
#    (red, differs from the key) as appropriate. Cells that stay blank
#    (not attempted) are left untouched.
# ---------------------------------------------------------------------------

# D16:D18 also hold a (correct) "Student Ans" duplicate that needs to be
# restyled/filled in the same way as column A.
$ws.Range("B10").Copy()
$ws.Range("D16:D18").PasteSpecial($xlPasteFormats)
$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Value = "Option D"

# Correct answers (green / correctStyle, same style already used by B10).
$correctCells = "A19", "A21", "A22", "A23", "A24", "A26", "A27", "A28", "A29", "A31", "A32", "A37", "A38", "A40"
$correctValues = @{
    "A19" = "Option C"
    "A21" = "Option C"
    "A22" = "Option D"
    "A23" = "Option D"
    "A24" = "Option A"
    "A26" = "Option C"
    "A27" = "Option A"
    "A28" = "Option D"
    "A29" = "Option D"
    "A31" = "Option D"
    "A32" = "Option C"
    "A37" = "Option A"
    "A38" = "Option A"
    "A40" = "Option D"
}
$ws.Range("B10").Copy()
foreach ($cellRef in $correctCells) {
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats)
}
foreach ($cellRef in $correctCells) {
    $ws.Range($cellRef).Value = $correctValues[$cellRef]
}

# Incorrect answers (red / incorrectStyle, same style already used by C10).
$incorrectValues = @{
    "A30" = "Option A"
    "A39" = "Option C"
}
$ws.Range("C10").Copy()
foreach ($cellRef in $incorrectValues.Keys) {
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats)
}
foreach ($cellRef in $incorrectValues.Keys) {
    $ws.Range($cellRef).Value = $incorrectValues[$cellRef]
}
